$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily-log row (row 5) under the existing battery_info data.
# Columns A,B,C,D,E,G,H hold plain text in this sheet (the date, time,
# percentage and cycle-count values are all stored as strings, not native
# Excel dates/numbers) while column F is a real number. A leading apostrophe
# forces Excel to keep date-/number-looking text ("2023-07-08", "19:34:21",
# "96%", "52") as literal text instead of auto-converting it; resetting the
# cell style back to "Normal" afterwards drops the quote-prefix formatting
# flag so the cell doesn't end up with a stray style index.
$ws.Range("A5").Value = "'2023-07-08"
$ws.Range("B5").Value = "'19:34:21"
$ws.Range("C5").Value = "Plugged"
$ws.Range("D5").Value = "charged;"
$ws.Range("E5").Value = "Normal"
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = "'96%"
$ws.Range("H5").Value = "'52"

$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Style = "Normal"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Style = "Normal"
